$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10 (columns E through T)
# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.29347
$ws.Range("H2").Value = 0.8804099999999999
$ws.Range("I2").Value = 0.1501202107524681
$ws.Range("J2").Value = 0.1501202107524681
$ws.Range("M2").Value = 0.7682129999999999
$ws.Range("N2").Value = 2.304639
$ws.Range("O2").Value = 0.06906161725690135
$ws.Range("P2").Value = 0.06906161725690135
$ws.Range("Q2").Value = 0.2254474691099999
$ws.Range("R2").Value = 2.02902722199
$ws.Range("S2").Value = 0.01036754453751232
$ws.Range("T2").Value = 0.01036754453751232

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.29347
$ws.Range("H3").Value = 0.8804099999999999
$ws.Range("I3").Value = 0.1501202107524681
$ws.Range("J3").Value = 0.1501202107524681
$ws.Range("M3").Value = 8.345897000000001
$ws.Range("O3").Value = 0.7502881938726906
$ws.Range("P3").Value = 0.7502881938726906
$ws.Range("Q3").Value = 2.44927039259
$ws.Range("R3").Value = 22.04343353331
$ws.Range("S3").Value = 0.112633421789257
$ws.Range("T3").Value = 0.112633421789257

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.29347
$ws.Range("H4").Value = 0.8804099999999999
$ws.Range("I4").Value = 0.1501202107524681
$ws.Range("J4").Value = 0.1501202107524681
$ws.Range("M4").Value = 2.009478333333333
$ws.Range("N4").Value = 6.028435
$ws.Range("O4").Value = 0.180650188870408
$ws.Range("P4").Value = 0.180650188870408
$ws.Range("Q4").Value = 0.5897216064833333
$ws.Range("R4").Value = 5.30749445835
$ws.Range("S4").Value = 0.02711924442569882
$ws.Range("T4").Value = 0.02711924442569882

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.2180391153852712
$ws.Range("J5").Value = 0.2180391153852712
$ws.Range("M5").Value = 0.7682129999999999
$ws.Range("N5").Value = 2.304639
$ws.Range("O5").Value = 0.06906161725690135
$ws.Range("P5").Value = 0.06906161725690135
$ws.Range("Q5").Value = 0.327446694114
$ws.Range("R5").Value = 2.947020247026
$ws.Range("S5").Value = 0.01505813393377095
$ws.Range("T5").Value = 0.01505813393377095

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.2180391153852712
$ws.Range("J6").Value = 0.2180391153852712
$ws.Range("M6").Value = 8.345897000000001
$ws.Range("O6").Value = 0.7502881938726906
$ws.Range("P6").Value = 0.7502881938726906
$ws.Range("Q6").Value = 3.557394084799334
$ws.Range("R6").Value = 32.01654676319401
$ws.Range("S6").Value = 0.1635921740760143
$ws.Range("T6").Value = 0.1635921740760143

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.2180391153852712
$ws.Range("J7").Value = 0.2180391153852712
$ws.Range("M7").Value = 2.009478333333333
$ws.Range("N7").Value = 6.028435
$ws.Range("O7").Value = 0.180650188870408
$ws.Range("P7").Value = 0.180650188870408
$ws.Range("Q7").Value = 0.8565294223655556
$ws.Range("R7").Value = 7.70876480129
$ws.Range("S7").Value = 0.03938880737548592
$ws.Range("T7").Value = 0.03938880737548592

# Row 8: MuSCs -> ECs
$ws.Range("G8").Value = 1.235185333333333
$ws.Range("H8").Value = 3.705556
$ws.Range("I8").Value = 0.6318406738622607
$ws.Range("J8").Value = 0.6318406738622606
$ws.Range("M8").Value = 0.7682129999999999
$ws.Range("N8").Value = 2.304639
$ws.Range("O8").Value = 0.06906161725690135
$ws.Range("P8").Value = 0.06906161725690135
$ws.Range("Q8").Value = 0.9488854304759999
$ws.Range("R8").Value = 8.539968874284
$ws.Range("S8").Value = 0.04363593878561808
$ws.Range("T8").Value = 0.04363593878561808

# Row 9: MuSCs -> FAPs
$ws.Range("G9").Value = 1.235185333333333
$ws.Range("H9").Value = 3.705556
$ws.Range("I9").Value = 0.6318406738622607
$ws.Range("J9").Value = 0.6318406738622606
$ws.Range("M9").Value = 8.345897000000001
$ws.Range("O9").Value = 0.7502881938726906
$ws.Range("P9").Value = 0.7502881938726906
$ws.Range("Q9").Value = 10.30872956791067
$ws.Range("R9").Value = 92.77856611119601
$ws.Range("S9").Value = 0.4740625980074193
$ws.Range("T9").Value = 0.4740625980074193

# Row 10: MuSCs -> MuSCs
$ws.Range("G10").Value = 1.235185333333333
$ws.Range("H10").Value = 3.705556
$ws.Range("I10").Value = 0.6318406738622607
$ws.Range("J10").Value = 0.6318406738622606
$ws.Range("M10").Value = 2.009478333333333
$ws.Range("N10").Value = 6.028435
$ws.Range("O10").Value = 0.180650188870408
$ws.Range("P10").Value = 0.180650188870408
$ws.Range("Q10").Value = 2.482078164984445
$ws.Range("R10").Value = 22.33870348486
$ws.Range("S10").Value = 0.1141421370692233
$ws.Range("T10").Value = 0.1141421370692232
